$d = $word.ActiveDocument

# Locate the paragraph that reads "Página de Requisitos FECHADA – data 16/05/2013"
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*Requisitos*FECHADA*") {
        $targetIndex = $i
        break
    }
}

$p = $d.Paragraphs.Item($targetIndex)
$r = $p.Range

# Bump the whole paragraph (including the paragraph mark) to 14pt (sz/szCs = 28)
$r.Font.Size = 14
$r.Font.SizeBi = 14

# Append a run of 13 trailing spaces after the existing content (after the
# _GoBack bookmark, before the paragraph mark).
$r.InsertAfter("             ")

# Re-fetch the paragraph and re-apply the 14pt size to the whole paragraph
# (including the just-inserted run) so the new trailing run also carries
# <w:sz>/<w:szCs>.
$p2 = $d.Paragraphs.Item($targetIndex)
$r2 = $p2.Range
$r2.Font.Size = 14
$r2.Font.SizeBi = 14
